# Auto-generated Excel COM-interop script
# Applies scheduled-runner price/profit updates to the Gungnir_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 485.85715
$ws.Range("I12").Value = 152.75
$ws.Range("J12").Value = 930
$ws.Range("K12").Value = 152.75
$ws.Range("L12").Value = 930
$ws.Range("M12").Value = 17.25
$ws.Range("N12").Value = -1270

$ws.Range("H19").Value = 500.18182
$ws.Range("I19").Value = 395.5
$ws.Range("J19").Value = 560
$ws.Range("K19").Value = 395.5
$ws.Range("L19").Value = 560
$ws.Range("M19").Value = -220.5
$ws.Range("N19").Value = -910

$ws.Range("H31").Value = 947.25
$ws.Range("J31").Value = 1000
$ws.Range("L31").Value = 3000
$ws.Range("N31").Value = -3460

$ws.Range("H40").Value = 626784.3
$ws.Range("I40").Value = 1371.8889
$ws.Range("J40").Value = 1430886
$ws.Range("K40").Value = 1371.8889
$ws.Range("L40").Value = 1430886
$ws.Range("M40").Value = -1196.8889
$ws.Range("N40").Value = -1431236

$ws.Range("H57").Value = 23695
$ws.Range("J57").Value = 23695
$ws.Range("L57").Value = 71085
$ws.Range("N57").Value = -72083

$ws.Range("H64").Value = 4000
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 4000
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 4000
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -4496

$ws.Range("H67").Value = 4000
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 4000
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 4000
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -5716

$ws.Range("H132").Value = 16138920
$ws.Range("I132").Value = 16676867
$ws.Range("J132").Value = 500
$ws.Range("K132").Value = 50030601
$ws.Range("L132").Value = 1500
$ws.Range("M132").Value = -50028071
$ws.Range("N132").Value = -6560

$ws.Range("H136").Value = 64800
$ws.Range("J136").Value = 64800
$ws.Range("L136").Value = 64800
$ws.Range("N136").Value = -75000

$ws.Range("H137").Value = 2251.3877
$ws.Range("I137").Value = 2212.5
$ws.Range("J137").Value = 2339.5334
$ws.Range("K137").Value = 6637.5
$ws.Range("L137").Value = 7018.600199999999
$ws.Range("M137").Value = -4087.5
$ws.Range("N137").Value = -12118.6002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H30").Value = 3333.1667
$ws.Range("I30").Value = 1666.6666
$ws.Range("J30").Value = 4999.6665
$ws.Range("K30").Value = 1666.6666
$ws.Range("L30").Value = 4999.6665
$ws.Range("M30").Value = -1516.6666
$ws.Range("N30").Value = -5299.6665

$ws.Range("H61").Value = 43270972
$ws.Range("I61").Value = 75001256
$ws.Range("J61").Value = 2395.4546
$ws.Range("K61").Value = 75001256
$ws.Range("L61").Value = 2395.4546
$ws.Range("M61").Value = -75001044
$ws.Range("N61").Value = -2819.4546

$ws.Range("H122").Value = 1439.9459
$ws.Range("I122").Value = 1478.7941
$ws.Range("J122").Value = 999.6667
$ws.Range("K122").Value = 4436.3823
$ws.Range("L122").Value = 2999.0001
$ws.Range("M122").Value = -1986.3823
$ws.Range("N122").Value = -7899.0001

$ws.Range("H136").Value = 43270972
$ws.Range("I136").Value = 75001256
$ws.Range("J136").Value = 2395.4546
$ws.Range("K136").Value = 225003768
$ws.Range("L136").Value = 7186.3638
$ws.Range("M136").Value = -225001218
$ws.Range("N136").Value = -12286.3638

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H109").Value = 30669
$ws.Range("J109").Value = 30669
$ws.Range("L109").Value = 30669
$ws.Range("N109").Value = -33443

$ws.Range("H117").Value = 30742
$ws.Range("J117").Value = 30742
$ws.Range("L117").Value = 30742
$ws.Range("N117").Value = -39920

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 587061.5
$ws.Range("I6").Value = 782415
$ws.Range("J6").Value = 1001
$ws.Range("K6").Value = 782415
$ws.Range("L6").Value = 1001
$ws.Range("M6").Value = -782302
$ws.Range("N6").Value = -1227

$ws.Range("H19").Value = 1959.75
$ws.Range("I19").Value = 1525.4286
$ws.Range("K19").Value = 1525.4286
$ws.Range("M19").Value = -1355.4286

$ws.Range("H24").Value = 1959.75
$ws.Range("I24").Value = 1525.4286
$ws.Range("K24").Value = 1525.4286
$ws.Range("M24").Value = -1355.4286

$ws.Range("H31").Value = 1070.9642
$ws.Range("I31").Value = 874.1
$ws.Range("J31").Value = 1180.3334
$ws.Range("K31").Value = 874.1
$ws.Range("L31").Value = 1180.3334
$ws.Range("M31").Value = -579.1
$ws.Range("N31").Value = -1770.3334

$ws.Range("H34").Value = 1070.9642
$ws.Range("I34").Value = 874.1
$ws.Range("J34").Value = 1180.3334
$ws.Range("K34").Value = 874.1
$ws.Range("L34").Value = 1180.3334
$ws.Range("M34").Value = -672.1
$ws.Range("N34").Value = -1584.3334

$ws.Range("H58").Value = 29412808
$ws.Range("I58").Value = 50000880
$ws.Range("J58").Value = 1276.9286
$ws.Range("K58").Value = 50000880
$ws.Range("L58").Value = 1276.9286
$ws.Range("M58").Value = -50000677
$ws.Range("N58").Value = -1682.9286

$ws.Range("H136").Value = 29412808
$ws.Range("I136").Value = 50000880
$ws.Range("J136").Value = 1276.9286
$ws.Range("K136").Value = 150002640
$ws.Range("L136").Value = 3830.7858
$ws.Range("M136").Value = -150000090
$ws.Range("N136").Value = -8930.7858

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 4113.8125
$ws.Range("I132").Value = 440.76923
$ws.Range("J132").Value = 6626.9473
$ws.Range("K132").Value = 3966.92307
$ws.Range("L132").Value = 59642.5257
$ws.Range("M132").Value = -1436.92307
$ws.Range("N132").Value = -64702.5257

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws.Range("H107").Value = 4493.4585
$ws.Range("J107").Value = 6594.5
$ws.Range("L107").Value = 6594.5
$ws.Range("N107").Value = -10434.5

$ws.Range("H122").Value = 22738482
$ws.Range("I122").Value = 41684784
$ws.Range("J122").Value = 2919.7
$ws.Range("K122").Value = 125054352
$ws.Range("L122").Value = 8759.099999999999
$ws.Range("M122").Value = -125051902
$ws.Range("N122").Value = -13659.1

$ws.Range("H123").Value = 10323.4
$ws.Range("J123").Value = 10323.4
$ws.Range("L123").Value = 10323.4
$ws.Range("N123").Value = -15223.4

$ws.Range("H132").Value = 5420.6587
$ws.Range("I132").Value = 2686.9644
$ws.Range("J132").Value = 11308.615
$ws.Range("K132").Value = 8060.8932
$ws.Range("L132").Value = 33925.845
$ws.Range("M132").Value = -5530.8932
$ws.Range("N132").Value = -38985.845

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1352.7646
$ws.Range("I68").Value = 1315.4375
$ws.Range("J68").Value = 1950
$ws.Range("K68").Value = 1315.4375
$ws.Range("L68").Value = 1950
$ws.Range("M68").Value = -566.4375
$ws.Range("N68").Value = -3448

$ws.Range("H71").Value = 1352.7646
$ws.Range("I71").Value = 1315.4375
$ws.Range("J71").Value = 1950
$ws.Range("K71").Value = 6577.1875
$ws.Range("L71").Value = 9750
$ws.Range("M71").Value = -2833.1875
$ws.Range("N71").Value = -17238

$ws.Range("H82").Value = 1242.8572
$ws.Range("I82").Value = 1260
$ws.Range("J82").Value = 1200
$ws.Range("K82").Value = 1260
$ws.Range("L82").Value = 1200
$ws.Range("M82").Value = -899
$ws.Range("N82").Value = -1922

$ws.Range("H85").Value = 1242.8572
$ws.Range("I85").Value = 1260
$ws.Range("J85").Value = 1200
$ws.Range("K85").Value = 1260
$ws.Range("L85").Value = 1200
$ws.Range("M85").Value = -12
$ws.Range("N85").Value = -3696

$ws.Range("H109").Value = 34910.4
$ws.Range("J109").Value = 34910.4
$ws.Range("L109").Value = 34910.4
$ws.Range("N109").Value = -37684.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 853.0769
$ws.Range("I81").Value = 588.6667
$ws.Range("K81").Value = 1177.3334
$ws.Range("M81").Value = -116.3334

$ws.Range("H84").Value = 853.0769
$ws.Range("I84").Value = 588.6667
$ws.Range("K84").Value = 5886.666999999999
$ws.Range("M84").Value = -582.6669999999995

$ws.Range("H132").Value = 29865.809
$ws.Range("I132").Value = 37155.934
$ws.Range("J132").Value = 11640.5
$ws.Range("K132").Value = 111467.802
$ws.Range("L132").Value = 34921.5
$ws.Range("M132").Value = -108937.802
$ws.Range("N132").Value = -39981.5

$ws.Range("H133").Value = 34238.332
$ws.Range("J133").Value = 34238.332
$ws.Range("L133").Value = 34238.332
$ws.Range("N133").Value = -44358.332
